# Update "想去人数" (F column) counts on both the "展览" and "全部类型"
# worksheets to match the regenerated data output.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 8818
    3  = 8170
    5  = 197
    8  = 141
    11 = 246
    12 = 738
    13 = 199
    14 = 4258
    16 = 77
    20 = 126
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
